# Fill in row 11 of Sheet1 with the new coverage-history data point,
# mirroring the pattern established by the preceding rows (e.g. row 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Carry the formatting from the row above down into row 11 (date format on
# column A, 2-decimal number format on the computed percentage columns)
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("N10:R10").Copy()
$ws.Range("N11:R11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Raw input values for the new week (row 11)
$ws.Range("A11").Value = 44082
$ws.Range("B11").Value = 10
$ws.Range("D11").Value = 162
$ws.Range("E11").Value = 84
$ws.Range("F11").Value = 1455
$ws.Range("G11").Value = 564
$ws.Range("H11").Value = 46
$ws.Range("I11").Value = 20
$ws.Range("J11").Value = 14
$ws.Range("K11").Value = 210
$ws.Range("L11").Value = 83

# Computed percentage columns, following the same formulas used in row 10
$ws.Range("N11").Formula = "=100*E11/D11"
$ws.Range("O11").Formula = "=100*G11/F11"
$ws.Range("P11").Formula = "=100*H11/D11"
$ws.Range("Q11").Formula = "=100*J11/I11"
$ws.Range("R11").Formula = "=100*L11/K11"

# The chart was nudged slightly up/left (one row and one column) by the
# author while reviewing the new data point. Re-anchor it from D7 to M39
# (it previously ran from E8 to N40), preserving its on-sheet size.
$co = $ws.ChartObjects().Item(1)
$fromCell = $ws.Cells.Item(7, 4)
$toCell = $ws.Cells.Item(39, 13)
$co.Left = $fromCell.Left + (247649 / 12700)
$co.Top = $fromCell.Top + (95250 / 12700)
$co.Width = ($toCell.Left + (209550 / 12700)) - $co.Left
$co.Height = ($toCell.Top + (19050 / 12700)) - $co.Top

# Leave the active selection where the author left it after entering the data
$ws.Range("M20").Select()

$wb.Save()
